$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (currently at the start of the
#    second paragraph, right before "The author of the article...").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# The paragraph containing "I tried to write something as usual, but I
# have failed. What can I say?" -- scope all edits to just this paragraph
# so we don't touch the other occurrences of "have " elsewhere in the doc.
$p3 = $d.Paragraphs.Item(3).Range

# 2. Delete the word "have " from "...but I have failed..." so the
#    sentence reads "...but I failed...".
$p3.Find.Execute("have ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 3. Re-insert the "_GoBack" bookmark right after "...but I " (i.e. right
#    before "failed."), splitting the trailing space into its own run.
$r = $d.Paragraphs.Item(3).Range
$r.Find.Execute("failed.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $r.Start

# Force the trailing space before "failed." to become its own run by
# toggling a character property off/on (net no-op) on just that character.
$spaceRange = $d.Range($pos - 1, $pos)
$spaceRange.Font.Bold = 1
$spaceRange.Font.Bold = 0

# Insert the bookmark collapsed right after that space run, before "failed.".
$insertRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $insertRange)
